# Updates cryptos list figures (price + 1h volume change) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "299.34")
    # are not silently coerced into floating point numbers, then
    # drop back to the default style so no stray NumberFormat sticks.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "23.447.56"
Set-TextValue "E2" "  -1.08%  "
Set-TextValue "D3" "1.646.73"
Set-TextValue "E3" "  -0.31%  "
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "E5" "  +0.22%  "
Set-TextValue "D6" "299.34"
Set-TextValue "E6" "  -1.54%  "
Set-TextValue "D7" "0.3792"
Set-TextValue "E7" "  -0.58%  "
Set-TextValue "D8" "0.3558"
Set-TextValue "E8" "  -1.36%  "
Set-TextValue "D9" "49.82"
Set-TextValue "E9" "  -3.25%  "
Set-TextValue "D10" "0.08100"
Set-TextValue "E10" "  -1.62%  "
Set-TextValue "D11" "1.219"
Set-TextValue "E11" "  -2.39%  "
Set-TextValue "D12" "1.003"
Set-TextValue "E12" "  +0.31%  "
Set-TextValue "D13" "22.02"
Set-TextValue "E13" "  -2.53%  "
Set-TextValue "D14" "6.388"
Set-TextValue "E14" "  -2.37%  "
Set-TextValue "D15" "7.361"
Set-TextValue "E15" "  -0.56%  "
Set-TextValue "D16" "0.00001194"
Set-TextValue "E16" "  -3.13%  "
Set-TextValue "D17" "1.642.83"
Set-TextValue "E17" "  -0.61%  "
Set-TextValue "D18" "97.33"
Set-TextValue "E18" "  +0.25%  "
Set-TextValue "D19" "0.06956"
Set-TextValue "E19" "  -0.17%  "
Set-TextValue "D20" "6.754"
Set-TextValue "E20" "  -0.52%  "
Set-TextValue "E21" "  -2.28%  "
Set-TextValue "D22" "1.002"
Set-TextValue "E22" "  +0.18%  "
Set-TextValue "D23" "12.38"
Set-TextValue "E23" "  -1.79%  "
Set-TextValue "D24" "23.451.85"
Set-TextValue "E24" "  -1.08%  "
Set-TextValue "D25" "2.498"
Set-TextValue "E25" "  -1.96%  "
Set-TextValue "D26" "2.936"
Set-TextValue "E26" "  -4.61%  "
Set-TextValue "D27" "20.89"
Set-TextValue "D28" "152.73"
Set-TextValue "E28" "  +0.30%  "
Set-TextValue "D29" "5.209"
Set-TextValue "E29" "  -0.52%  "
Set-TextValue "D30" "132.80"
Set-TextValue "E30" "  -1.74%  "
Set-TextValue "D31" "1.837.43"
Set-TextValue "E31" "  +0.01%  "
Set-TextValue "D32" "6.911"
Set-TextValue "E32" "  +0.44%  "
Set-TextValue "D33" "2.114"
Set-TextValue "E33" "  +0.53%  "
Set-TextValue "D34" "11.77"
Set-TextValue "E34" "  -1.66%  "
Set-TextValue "D35" "1.010"
Set-TextValue "E35" "  -6.95%  "
Set-TextValue "D36" "0.02720"
Set-TextValue "E36" "  -3.50%  "
Set-TextValue "D37" "0.08734"
Set-TextValue "E37" "  -1.20%  "
Set-TextValue "D38" "0.2428"
Set-TextValue "E38" "  -3.60%  "
Set-TextValue "D39" "5.927"
Set-TextValue "E39" "  -2.68%  "
Set-TextValue "D40" "13.09"
Set-TextValue "E40" "  +1.95%  "
Set-TextValue "D41" "0.06784"
Set-TextValue "E41" "  -3.79%  "
Set-TextValue "D42" "0.6870"
Set-TextValue "E42" "  -2.68%  "
Set-TextValue "D43" "1.301"
Set-TextValue "E43" "  -2.89%  "
Set-TextValue "D44" "15.58"
Set-TextValue "E44" "  -2.53%  "
Set-TextValue "E45" "  +0.23%  "
Set-TextValue "D46" "0.6393"
Set-TextValue "E46" "  -1.90%  "
Set-TextValue "D47" "2.250"
Set-TextValue "E47" "  -3.83%  "
Set-TextValue "D48" "3.916"
Set-TextValue "E48" "  -1.69%  "
Set-TextValue "D49" "0.07717"
Set-TextValue "E49" "  -3.39%  "
Set-TextValue "D50" "127.31"
Set-TextValue "E50" "  -0.63%  "
Set-TextValue "D51" "1.149"
Set-TextValue "E51" "  -3.40%  "
